$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at M; this shifts the existing M (total_files) to N
# and the existing N (unique_days) to O, and keeps their widths/styles intact.
$ws.Columns("M").Insert()

# Restore the explicit column width for the freshly inserted column M
# (Excel's insert leaves it at the default width).
$ws.Columns("M").ColumnWidth = 11.2

# New header cell for the 2026-02-27 submissions column. The leading
# apostrophe forces the date-shaped text to be stored as text (matching
# the other date header cells) instead of being auto-parsed into a date
# serial number.
$ws.Range("M1").Value = "'2026-02-27"

# Per-student data: new 2026-02-27 indicator (M), recomputed total_files
# (N = old total_files + new indicator) and recomputed unique_days
# (O = old unique_days + (1 if they submitted on 2026-02-27 else 0)).
$data = @(
    @(2,1,6,6),
    @(3,1,5,5),
    @(4,1,6,6),
    @(5,1,7,7),
    @(6,1,5,5),
    @(7,1,7,7),
    @(8,1,7,7),
    @(9,1,7,7),
    @(10,1,6,6),
    @(11,1,6,6),
    @(12,1,7,7),
    @(13,0,5,5),
    @(14,1,7,7),
    @(15,1,5,5),
    @(16,1,5,5),
    @(17,1,7,7),
    @(18,1,7,7),
    @(19,1,5,5),
    @(20,1,6,6),
    @(21,0,0,0),
    @(22,0,0,0),
    @(23,1,2,2),
    @(24,1,6,6),
    @(25,1,5,5),
    @(26,1,4,4),
    @(27,1,8,6),
    @(28,0,0,0),
    @(29,0,0,0),
    @(30,0,0,0),
    @(31,1,7,7),
    @(32,1,7,7),
    @(33,1,7,7),
    @(34,0,6,6),
    @(35,0,5,5),
    @(36,1,6,6),
    @(37,0,3,3),
    @(38,0,0,0),
    @(39,0,0,0),
    @(40,1,7,7),
    @(41,1,7,7),
    @(42,1,7,7),
    @(43,0,21,1),
    @(44,1,45,4),
    @(45,1,4,4),
    @(46,1,7,7),
    @(47,0,6,6),
    @(48,0,0,0),
    @(49,0,6,6),
    @(50,1,6,6),
    @(51,1,2,2),
    @(52,1,3,3),
    @(53,1,6,6),
    @(54,0,0,0),
    @(55,0,6,6),
    @(56,0,0,0),
    @(57,1,6,6),
    @(58,1,6,6),
    @(59,0,3,3),
    @(60,1,7,7),
    @(61,0,1,1),
    @(62,0,0,0),
    @(63,1,4,4),
    @(64,1,7,7),
    @(65,0,0,0),
    @(66,0,0,0),
    @(67,1,3,3),
    @(68,0,0,0),
    @(69,0,0,0),
    @(70,1,7,7),
    @(71,0,1,1),
    @(72,0,0,0),
    @(73,0,3,3),
    @(74,1,18,6),
    @(75,1,2,2),
    @(76,1,5,5),
    @(77,0,0,0),
    @(78,0,0,0),
    @(79,0,3,3),
    @(80,1,6,6),
    @(81,1,5,5),
    @(82,1,7,7),
    @(83,1,4,4),
    @(84,0,1,1),
    @(85,0,0,0),
    @(86,1,4,4),
    @(87,0,0,0),
    @(88,1,3,3),
    @(89,0,0,0),
    @(90,0,2,2),
    @(91,0,0,0),
    @(92,0,0,0),
    @(93,1,4,4),
    @(94,1,5,5),
    @(95,1,34,5),
    @(96,1,2,2),
    @(97,0,0,0),
    @(98,1,3,3),
    @(99,0,0,0),
    @(100,0,2,2),
    @(101,1,7,7),
    @(102,0,0,0),
    @(103,1,7,7),
    @(104,0,12,1),
    @(105,1,7,7),
    @(106,0,6,6),
    @(107,0,0,0),
    @(108,1,6,6),
    @(109,0,0,0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 13).Value = $row[1]
    $ws.Cells.Item($r, 14).Value = $row[2]
    $ws.Cells.Item($r, 15).Value = $row[3]
}
